$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# --- Copy number-format/style from H:J onto K:M so the style index matches (xf idx 6 -> 8) ---
foreach ($r in 9..25) {
    $ws.Range("H" + $r + ":J" + $r).Copy()
    $ws.Range("K" + $r + ":M" + $r).PasteSpecial(-4122)
}
$ws.Range("H27:J27").Copy()
$ws.Range("K27:M27").PasteSpecial(-4122)
$ws.Range("J26").Copy()
$ws.Range("M26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Write the updated values (rolling-window recompute) ---
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = -1

$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = -1
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = -0.94736842

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = -0.94736842
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = -0.89473684

$ws.Range("H12").Value = 0
$ws.Range("J12").Value = -0.89473684
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = -0.84210526

$ws.Range("J13").Value = -0.84210526
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = -0.78947368

$ws.Range("J14").Value = -0.78947368
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = -0.73684211

$ws.Range("J15").Value = -0.73684211
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = -0.68421053

$ws.Range("J16").Value = -0.68421053
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = -0.63157895

$ws.Range("J17").Value = -0.63157895
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = -0.57894737

$ws.Range("J18").Value = -0.57894737
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = -0.52631579

$ws.Range("J19").Value = -0.52631579
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = -0.47368421

$ws.Range("J20").Value = -0.47368421
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = -0.42105263

$ws.Range("J21").Value = -0.42105263
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = -0.36842105

$ws.Range("J22").Value = -0.36842105
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = -0.31578947

$ws.Range("J23").Value = -0.31578947
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = -0.26315789

$ws.Range("J24").Value = -0.26315789
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = -0.21052632

$ws.Range("J25").Value = -0.21052632
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = -0.15789474

$ws.Range("J26").Value = -0.15789474
$ws.Range("M26").Value = -0.10526315

$ws.Range("J27").Value = -0.10526315
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = -0.05263156

# --- Match the recorded selection at save time ---
$ws.Range("M27").Select()

